$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first column (A) contained redundant row-index values (1, 17) that
# duplicated the data already present in column F. Delete column A entirely
# so every other column shifts one position to the left.
$ws.Range("A:A").Delete()
